$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'243.55"
$ws.Cells.Item(2, 4).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 4).Value = "'23.89"
$ws.Cells.Item(3, 4).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 2).Value = "HuobiToken"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(4, 4).Value = "'5.259"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "3HuobiTokenHT"

# Row 5
$ws.Cells.Item(5, 2).Value = "Cronos"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(5, 4).Value = "'0.05820"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "4CronosCRO"

# Row 6
$ws.Cells.Item(6, 2).Value = "KuCoinToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(6, 4).Value = "'6.460"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "5KuCoinTokenKCS"

# Row 7
$ws.Cells.Item(7, 2).Value = "GateToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7, 4).Value = "'3.332"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "6GateTokenGT"

# Row 8
$ws.Cells.Item(8, 2).Value = "MXToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8, 4).Value = "'0.8085"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "7MXTokenMX"

# Row 9
$ws.Cells.Item(9, 2).Value = "FTXToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(9, 4).Value = "'0.8738"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "8FTXTokenFTT"

# Row 10
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Value = "'0.1382"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "9WazirXWRX"

# Row 11
$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).Value = "'0.07270"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "10MandalaExchangeTokenMDX"

# Row 12
$ws.Cells.Item(12, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(12, 4).Value = "'0.03064"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).Value = "'0.03051"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

# Row 14
$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 4).Value = "'0.09331"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"

# Row 15
$ws.Cells.Item(15, 2).Value = "MCDex"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(15, 4).Value = "'3.853"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "14MCDexMCB"

# Row 16
$ws.Cells.Item(16, 2).Value = "BitForexToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(16, 4).Value = "'0.001535"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "15BitForexTokenBF"

# Row 17
$ws.Cells.Item(17, 2).Value = "CoinExToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(17, 4).Value = "'0.04717"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "16CoinExTokenCET"

# Row 18
$ws.Cells.Item(18, 2).Value = "One"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(18, 4).Value = "'0.0006038"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "17OneONE"

# Row 19
$ws.Cells.Item(19, 4).Value = "'0.006149"
$ws.Cells.Item(19, 4).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 4).Value = "'0.001270"
$ws.Cells.Item(20, 4).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 4).Value = "'0.004592"
$ws.Cells.Item(21, 4).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.00008695"
$ws.Cells.Item(22, 4).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 2).Value = "LEO"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(23, 4).Value = "'3.562"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "22LEOLEO"

# Row 24
$ws.Cells.Item(24, 2).Value = "BTSEToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(24, 4).Value = "'2.184"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "23BTSETokenBTSE"

# Row 25
$ws.Cells.Item(25, 4).Value = "'0.3176"
$ws.Cells.Item(25, 4).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.1318"
$ws.Cells.Item(26, 4).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.03774"
$ws.Cells.Item(40, 4).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 4).Value = "'0.1055"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "40BKEXTokenBKK"

# Row 42
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42, 4).Value = "'0.002708"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"

# Row 43
$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43, 4).Value = "'0.003233"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "42KickTokenKICKWorstin24h"

# Row 44
$ws.Cells.Item(44, 4).Value = "'0.007973"
$ws.Cells.Item(44, 4).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.00005502"
$ws.Cells.Item(45, 4).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOIN"
